# Apply the "adaptive bot analyses" figure fix described in the commit:
#   5. Win-positive-lose-negative  ->  5. Win-positive lose-negative
#   6. Win-stay-lose-positive      ->  6. Win-stay lose-positive
#
# Both labels live as single text runs inside the "TextBox 10" shape on
# slide 1 (the legend explaining the outcome/opponent transition codes).
# We replace each run's text in place (same length, just one hyphen
# becomes a space) so every other run/paragraph/formatting attribute in
# the shape is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the legend textbox by name rather than a hard-coded index.
$legend = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 10") {
        $legend = $candidate
    }
}

$tr = $legend.TextFrame.TextRange

function Fix-HyphenToSpace($textRange, $oldText, $newText) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -ge 0) {
        $run = $textRange.Characters($idx + 1, $oldText.Length)
        $run.Text = $newText
    }
}

Fix-HyphenToSpace $tr "5. Win-positive-lose-negative" "5. Win-positive lose-negative"
Fix-HyphenToSpace $tr "6. Win-stay-lose-positive" "6. Win-stay lose-positive"
